$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 33.36960033333333
$ws.Range("H2").Value = 100.108801
$ws.Range("I2").Value = 0.07727383968381614
$ws.Range("J2").Value = 0.07727383968381614
$ws.Range("M2").Value = 1.619868333333333
$ws.Range("N2").Value = 4.859605
$ws.Range("O2").Value = 0.1089327058120143
$ws.Range("P2").Value = 0.1089327058120143
$ws.Range("Q2").Value = 54.05435887595611
$ws.Range("R2").Value = 486.489229883605
$ws.Range("S2").Value = 0.0084176484452419
$ws.Range("T2").Value = 0.008417648445241899
$ws.Range("G3").Value = 33.36960033333333
$ws.Range("H3").Value = 100.108801
$ws.Range("I3").Value = 0.07727383968381614
$ws.Range("J3").Value = 0.07727383968381614
$ws.Range("O3").Value = 0.1655705935257241
$ws.Range("P3").Value = 0.1655705935257241
$ws.Range("Q3").Value = 82.15909276309799
$ws.Range("R3").Value = 739.431834867882
$ws.Range("S3").Value = 0.01279427550046109
$ws.Range("T3").Value = 0.01279427550046109
$ws.Range("G4").Value = 33.36960033333333
$ws.Range("H4").Value = 100.108801
$ws.Range("I4").Value = 0.07727383968381614
$ws.Range("J4").Value = 0.07727383968381614
$ws.Range("M4").Value = 7.682722666666667
$ws.Range("N4").Value = 23.048168
$ws.Range("O4").Value = 0.5166467859527435
$ws.Range("P4").Value = 0.5166467859527435
$ws.Range("Q4").Value = 256.3693848585075
$ws.Range("R4").Value = 2307.324463726568
$ws.Range("S4").Value = 0.03992328091087118
$ws.Range("T4").Value = 0.03992328091087118
$ws.Range("G5").Value = 33.36960033333333
$ws.Range("H5").Value = 100.108801
$ws.Range("I5").Value = 0.07727383968381614
$ws.Range("J5").Value = 0.07727383968381614
$ws.Range("M5").Value = 3.105673
$ws.Range("N5").Value = 9.317019
$ws.Range("O5").Value = 0.2088499147095181
$ws.Range("P5").Value = 0.2088499147095181
$ws.Range("Q5").Value = 103.6350667760243
$ws.Range("R5").Value = 932.715600984219
$ws.Range("S5").Value = 0.01613863482724198
$ws.Range("T5").Value = 0.01613863482724198
$ws.Range("I6").Value = 0.2551852590901843
$ws.Range("J6").Value = 0.2551852590901843
$ws.Range("M6").Value = 1.619868333333333
$ws.Range("N6").Value = 4.859605
$ws.Range("O6").Value = 0.1089327058120143
$ws.Range("P6").Value = 0.1089327058120143
$ws.Range("Q6").Value = 178.5064082638511
$ws.Range("R6").Value = 1606.55767437466
$ws.Range("S6").Value = 0.0277980207560337
$ws.Range("T6").Value = 0.02779802075603369
$ws.Range("I7").Value = 0.2551852590901843
$ws.Range("J7").Value = 0.2551852590901843
$ws.Range("O7").Value = 0.1655705935257241
$ws.Range("P7").Value = 0.1655705935257241
$ws.Range("S7").Value = 0.04225117480657751
$ws.Range("T7").Value = 0.0422511748065775
$ws.Range("I8").Value = 0.2551852590901843
$ws.Range("J8").Value = 0.2551852590901843
$ws.Range("M8").Value = 7.682722666666667
$ws.Range("N8").Value = 23.048168
$ws.Range("O8").Value = 0.5166467859527435
$ws.Range("P8").Value = 0.5166467859527435
$ws.Range("Q8").Value = 846.6214202063396
$ws.Range("R8").Value = 7619.592781857056
$ws.Range("S8").Value = 0.1318406439314618
$ws.Range("T8").Value = 0.1318406439314618
$ws.Range("I9").Value = 0.2551852590901843
$ws.Range("J9").Value = 0.2551852590901843
$ws.Range("M9").Value = 3.105673
$ws.Range("N9").Value = 9.317019
$ws.Range("O9").Value = 0.2088499147095181
$ws.Range("P9").Value = 0.2088499147095181
$ws.Range("Q9").Value = 342.2392555395053
$ws.Range("R9").Value = 3080.153299855548
$ws.Range("S9").Value = 0.05329541959611127
$ws.Range("T9").Value = 0.05329541959611127
$ws.Range("G10").Value = 13.90116633333333
$ws.Range("H10").Value = 41.703499
$ws.Range("I10").Value = 0.0321908709702775
$ws.Range("J10").Value = 0.0321908709702775
$ws.Range("M10").Value = 1.619868333333333
$ws.Range("N10").Value = 4.859605
$ws.Range("O10").Value = 0.1089327058120143
$ws.Range("P10").Value = 0.1089327058120143
$ws.Range("Q10").Value = 22.51805913976612
$ws.Range("R10").Value = 202.662532257895
$ws.Range("S10").Value = 0.003506638677237751
$ws.Range("T10").Value = 0.00350663867723775
$ws.Range("G11").Value = 13.90116633333333
$ws.Range("H11").Value = 41.703499
$ws.Range("I11").Value = 0.0321908709702775
$ws.Range("J11").Value = 0.0321908709702775
$ws.Range("O11").Value = 0.1655705935257241
$ws.Range("P11").Value = 0.1655705935257241
$ws.Range("Q11").Value = 34.225978222302
$ws.Range("R11").Value = 308.033804000718
$ws.Range("S11").Value = 0.005329861612658849
$ws.Range("T11").Value = 0.005329861612658848
$ws.Range("G12").Value = 13.90116633333333
$ws.Range("H12").Value = 41.703499
$ws.Range("I12").Value = 0.0321908709702775
$ws.Range("J12").Value = 0.0321908709702775
$ws.Range("M12").Value = 7.682722666666667
$ws.Range("N12").Value = 23.048168
$ws.Range("O12").Value = 0.5166467859527435
$ws.Range("P12").Value = 0.5166467859527435
$ws.Range("Q12").Value = 106.7988056822036
$ws.Range("R12").Value = 961.189251139832
$ws.Range("S12").Value = 0.01663131002381334
$ws.Range("T12").Value = 0.01663131002381334
$ws.Range("G13").Value = 13.90116633333333
$ws.Range("H13").Value = 41.703499
$ws.Range("I13").Value = 0.0321908709702775
$ws.Range("J13").Value = 0.0321908709702775
$ws.Range("M13").Value = 3.105673
$ws.Range("N13").Value = 9.317019
$ws.Range("O13").Value = 0.2088499147095181
$ws.Range("P13").Value = 0.2088499147095181
$ws.Range("Q13").Value = 43.17247694994234
$ws.Range("R13").Value = 388.552292549481
$ws.Range("S13").Value = 0.006723060656567558
$ws.Range("T13").Value = 0.006723060656567558
$ws.Range("G14").Value = 274.366806
$ws.Range("H14").Value = 823.100418
$ws.Range("I14").Value = 0.635350030255722
$ws.Range("J14").Value = 0.635350030255722
$ws.Range("M14").Value = 1.619868333333333
$ws.Range("N14").Value = 4.859605
$ws.Range("O14").Value = 0.1089327058120143
$ws.Range("P14").Value = 0.1089327058120143
$ws.Range("Q14").Value = 444.43810075721
$ws.Range("R14").Value = 3999.94290681489
$ws.Range("S14").Value = 0.06921039793350095
$ws.Range("T14").Value = 0.06921039793350095
$ws.Range("G15").Value = 274.366806
$ws.Range("H15").Value = 823.100418
$ws.Range("I15").Value = 0.635350030255722
$ws.Range("J15").Value = 0.635350030255722
$ws.Range("O15").Value = 0.1655705935257241
$ws.Range("P15").Value = 0.1655705935257241
$ws.Range("Q15").Value = 675.516866851764
$ws.Range("R15").Value = 6079.651801665877
$ws.Range("S15").Value = 0.1051952816060267
$ws.Range("T15").Value = 0.1051952816060267
$ws.Range("G16").Value = 274.366806
$ws.Range("H16").Value = 823.100418
$ws.Range("I16").Value = 0.635350030255722
$ws.Range("J16").Value = 0.635350030255722
$ws.Range("M16").Value = 7.682722666666667
$ws.Range("N16").Value = 23.048168
$ws.Range("O16").Value = 0.5166467859527435
$ws.Range("P16").Value = 0.5166467859527435
$ws.Range("Q16").Value = 2107.884079437136
$ws.Range("R16").Value = 18970.95671493423
$ws.Range("S16").Value = 0.3282515510865971
$ws.Range("T16").Value = 0.3282515510865971
$ws.Range("G17").Value = 274.366806
$ws.Range("H17").Value = 823.100418
$ws.Range("I17").Value = 0.635350030255722
$ws.Range("J17").Value = 0.635350030255722
$ws.Range("M17").Value = 3.105673
$ws.Range("N17").Value = 9.317019
$ws.Range("O17").Value = 0.2088499147095181
$ws.Range("P17").Value = 0.2088499147095181
$ws.Range("Q17").Value = 852.0935814904379
$ws.Range("R17").Value = 7668.842233413942
$ws.Range("S17").Value = 0.1326927996295973
$ws.Range("T17").Value = 0.1326927996295973
